# Add the new "2022-Q3" quarterly sheet and record its summary row on "总计".
#
# Helper: force a cell to hold a literal TEXT value (never auto-coerced to a
# number) while leaving its style at the workbook default (no style index),
# matching how the source data was authored (values such as "001417" or
# "32.06" must stay text, not become 1417 / 32.06 numbers).
function Set-TextCell {
    param($Cell, [string]$Text, $BlankCell)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $BlankCell.Copy()
    $Cell.PasteSpecial(-4122)  # xlPasteFormats - restores default (no) style, keeps the text value
}

$wb = $excel.ActiveWorkbook
$zongji = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new worksheet "2022-Q3" right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $zongji)
$q3.Name = "2022-Q3"

$q3.PageSetup.LeftMargin = 0.75 * 72
$q3.PageSetup.RightMargin = 0.75 * 72
$q3.PageSetup.TopMargin = 1 * 72
$q3.PageSetup.BottomMargin = 1 * 72
$q3.PageSetup.HeaderMargin = 0.5 * 72
$q3.PageSetup.FooterMargin = 0.5 * 72

$blank = $q3.Cells.Item(5000, 200)

# Header row (style copied from the matching header on the "总计" sheet).
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$zongji.Cells.Item(1,2).Copy()
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).PasteSpecial(-4122)
    Set-TextCell $q3.Cells.Item(1, $col) $headers[$col - 2] $blank
    $zongji.Cells.Item(1,2).Copy()
    $q3.Cells.Item(1, $col).PasteSpecial(-4122)
}

# Data rows: code, name, scale, position, ratio, marketValue, rank
$rows = @(
    @("001417","汇添富医疗服务灵活配置混合A","32.06","85.81","3.44","1.1029","7"),
    @("000727","融通健康产业灵活配置混合A","22.64","93.67","2.90","0.6566","10"),
    @("009274","融通健康产业灵活配置混合C","17.64","93.67","2.90","0.5116","10"),
    @("012358","汇丰晋信医疗先锋混合A","1.87","58.14","4.04","0.0755","1"),
    @("002863","金信深圳成长灵活配置混合","0.61","93.95","9.80","0.0598","2"),
    @("011765","兴银高端制造混合A","0.57","92.99","4.66","0.0266","1"),
    @("540007","汇丰晋信中小盘股票","0.56","92.05","3.94","0.0221","5"),
    @("013441","西藏东财创新医疗六个月定开混合","0.49","82.53","3.97","0.0195","9"),
    @("011766","兴银高端制造混合C","0.34","92.99","4.66","0.0158","1"),
    @("003513","中邮消费升级灵活配置混合","0.56","30.56","2.34","0.0131","7"),
    @("001537","中加改革红利灵活配置混合","0.37","79.28","3.29","0.0122","8"),
    @("012359","汇丰晋信医疗先锋混合C","0.12","58.14","4.04","0.0048","1"),
    @("008037","兴银先锋成长混合A","0.21","71.76","2.07","0.0043","7"),
    @("008038","兴银先锋成长混合C","0.13","71.76","2.07","0.0027","7"),
    @("005146","兴银丰润灵活配置混合","0.04","92.81","4.52","0.0018","2"),
    @("015121","汇添富医疗服务灵活配置混合C","0.02","85.81","3.44","0.0007","7"),
    @("015122","汇添富医疗服务灵活配置混合D","0.00","85.81","3.44","0","7")
)

$zongji.Cells.Item(2,1).Copy()
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q3.Cells.Item($r, 1).PasteSpecial(-4122)
    $q3.Cells.Item($r, 1).Value = $i

    Set-TextCell $q3.Cells.Item($r, 2) $row[0] $blank
    Set-TextCell $q3.Cells.Item($r, 3) $row[1] $blank
    Set-TextCell $q3.Cells.Item($r, 4) $row[2] $blank
    Set-TextCell $q3.Cells.Item($r, 5) $row[3] $blank
    Set-TextCell $q3.Cells.Item($r, 6) $row[4] $blank

    if ($r -eq 18) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextCell $q3.Cells.Item($r, 7) $row[5] $blank
    }

    $q3.Cells.Item($r, 8).Value = [int]$row[6]
}

$blank.Clear()

# ---------------------------------------------------------------------
# 2. Record the 2022-Q3 summary row on "总计" (it becomes the new top row,
#    the rest of the quarters shift down by one).
# ---------------------------------------------------------------------
$summary = @(
    @("2022-Q3", 17, 2.53),
    @("2022-Q2", 32, 6.3),
    @("2022-Q1", 39, 19.55),
    @("2021-Q4", 68, 31.08),
    @("2021-Q3", 66, 29.9),
    @("2021-Q2", 38, 12.69),
    @("2021-Q1", 6, 0.49)
)

$zongji.Cells.Item(2,1).Copy()
$zongji.Cells.Item(8,1).PasteSpecial(-4122)

for ($i = 0; $i -lt $summary.Count; $i++) {
    $r = $i + 2
    $item = $summary[$i]
    $zongji.Cells.Item($r, 1).Value = $i
    $zongji.Cells.Item($r, 2).Value = $item[0]
    $zongji.Cells.Item($r, 3).Value = $item[1]
    $zongji.Cells.Item($r, 4).Value = $item[2]
}

$zongji.Activate()
